$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.029.78"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.240.68"
$ws.Range("E3").Value = "  -1.80%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'315.22"
$ws.Range("E5").Value = "  -1.70%  "
$ws.Range("D6").Value = "'99.28"
$ws.Range("E6").Value = "  -7.51%  "
$ws.Range("E7").Value = "  -3.19%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -7.49%  "
$ws.Range("D10").Value = "'36.15"
$ws.Range("E10").Value = "  -7.24%  "
$ws.Range("E11").Value = "  -2.67%  "
$ws.Range("E12").Value = "  -7.46%  "
$ws.Range("E13").Value = "  -2.99%  "
$ws.Range("D14").Value = "2.582.36"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "'0.842"
$ws.Range("E15").Value = "  -5.29%  "
$ws.Range("D16").Value = "2.235.87"
$ws.Range("E16").Value = "  -2.34%  "
$ws.Range("D17").Value = "'13.95"
$ws.Range("E17").Value = "  -5.08%  "
$ws.Range("D18").Value = "43.883.70"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").Value = "'13.18"
$ws.Range("E19").Value = "  -7.13%  "
$ws.Range("E20").Value = "  -3.20%  "
$ws.Range("E21").Value = "  -4.05%  "
$ws.Range("D22").Value = "'65.93"
$ws.Range("E22").Value = "  -1.20%  "
$ws.Range("D23").Value = "'237.64"
$ws.Range("E23").Value = "  -0.79%  "
$ws.Range("E24").Value = "  -8.12%  "
$ws.Range("E25").Value = "  -8.64%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'10.16"
$ws.Range("E27").Value = "  -1.18%  "
$ws.Range("D28").Value = "'2.13"
$ws.Range("E28").Value = "  -4.56%  "
$ws.Range("D29").Value = "'36.53"
$ws.Range("E29").Value = "  -6.21%  "
$ws.Range("E30").Value = "  -8.93%  "
$ws.Range("E31").Value = "  -3.58%  "
$ws.Range("D32").Value = "'154.25"
$ws.Range("E32").Value = "  -5.60%  "
$ws.Range("D33").Value = "'0.0834"
$ws.Range("E33").Value = "  -6.09%  "
$ws.Range("E34").Value = "  +2.20%  "
$ws.Range("E35").Value = "  -2.80%  "
$ws.Range("E36").Value = "  -8.00%  "
$ws.Range("E37").Value = "  -7.24%  "
$ws.Range("D38").Value = "'0.117"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "'15.82"
$ws.Range("E39").Value = "  +2.40%  "
$ws.Range("E40").Value = "  -11.81%  "
$ws.Range("E41").Value = "  -10.78%  "
$ws.Range("E42").Value = "  -6.74%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").Value = "1.704.60"
$ws.Range("E44").Value = "  -4.34%  "
$ws.Range("D45").Value = "'82.00"
$ws.Range("E45").Value = "  -5.30%  "
$ws.Range("E46").Value = "  -6.61%  "
$ws.Range("D47").Value = "'5.18"
$ws.Range("E47").Value = "  -5.56%  "
$ws.Range("D48").Value = "'101.71"
$ws.Range("E48").Value = "  -3.05%  "
$ws.Range("D49").Value = "'71.65"
$ws.Range("E49").Value = "  -5.27%  "
$ws.Range("D50").Value = "'56.38"
$ws.Range("E50").Value = "  -7.04%  "
$ws.Range("D51").Value = "'1.62"
$ws.Range("E51").Value = "  -4.92%  "
